# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Every player row gets the same team season record: 69 wins, 93 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
# Copy the formatting of the last existing header cell (AC1) onto the three
# new header cells so they pick up the same bold/centered/bordered style,
# then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-58) ----------------------------------------------
$ws.Range("AD2:AD58").Value = 69
$ws.Range("AE2:AE58").Value = 93
$ws.Range("AF2:AF58").Value = 0

Write-Output "done"
